$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before U. This shifts the old "param_P_to_charging_station1"
# column (U) to V, and the old "param_P_to_charging_station2" column (V) to W,
# carrying over the header style (bold/bordered) used by the rest of row 1.
$ws.Columns("U").Insert()

# New header for the inserted column U (shared string "param_E_pv3_solar")
$ws.Range("U1").Value = "param_E_pv3_solar"

# Fill data rows 2-17 for columns U (new param_E_pv3_solar series, constant 0.12),
# and the recalculated values for V (param_P_to_charging_station1) and
# W (param_P_to_charging_station2).
$ws.Range("U2").Value = 0.12
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("U3").Value = 0.12
$ws.Range("V3").Value = 91.43994444444448
$ws.Range("W3").Value = 23.15744444444445
$ws.Range("U4").Value = 0.12
$ws.Range("V4").Value = 89.48500000000003
$ws.Range("W4").Value = 122.4135000000001
$ws.Range("U5").Value = 0.12
$ws.Range("V5").Value = 38.11500000000002
$ws.Range("W5").Value = 0
$ws.Range("U6").Value = 0.12
$ws.Range("V6").Value = 20.66350000000001
$ws.Range("W6").Value = 0
$ws.Range("U7").Value = 0.12
$ws.Range("V7").Value = 42.35
$ws.Range("W7").Value = 52.40552777777779
$ws.Range("U8").Value = 0.12
$ws.Range("V8").Value = 69.46500000000002
$ws.Range("W8").Value = 65.41700000000003
$ws.Range("U9").Value = 0.12
$ws.Range("V9").Value = 66.00000000000004
$ws.Range("W9").Value = 34.73616666666666
$ws.Range("U10").Value = 0.12
$ws.Range("V10").Value = 0
$ws.Range("W10").Value = 0
$ws.Range("U11").Value = 0.12
$ws.Range("V11").Value = 0
$ws.Range("W11").Value = 111.9891666666667
$ws.Range("U12").Value = 0.12
$ws.Range("V12").Value = 49.93404166666668
$ws.Range("W12").Value = 33
$ws.Range("U13").Value = 0.12
$ws.Range("V13").Value = 33.37400000000001
$ws.Range("W13").Value = 0
$ws.Range("U14").Value = 0.12
$ws.Range("V14").Value = 157.15425
$ws.Range("W14").Value = 0
$ws.Range("U15").Value = 0.12
$ws.Range("V15").Value = 0
$ws.Range("W15").Value = 0
$ws.Range("U16").Value = 0.12
$ws.Range("V16").Value = 0
$ws.Range("W16").Value = 0
$ws.Range("U17").Value = 0.12
$ws.Range("V17").Value = 0
$ws.Range("W17").Value = 0

Write-Output "done"
